# Incorporación de nombres elecciones a clase Tablero
# Insert a new row (Gobernatura 17 / gb_17 / #dda15e) right after the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2, shifting existing data rows down.
$ws.Rows.Item(2).Insert()

# Populate the new row with the "Gobernatura 17" data.
$ws.Range("A2").Value = "Gobernatura 17"
$ws.Range("B2").Value = "gb_17"
$ws.Range("C2").Value = "#dda15e"

# Move the active selection to C2, matching the post-edit selection in the file.
$ws.Range("C2").Select()
